{"js": "/*\n * Office.js (Word JavaScript API) script implementing the resume rewrite:\n *  - Strip ad-hoc run colors/sizes from headings and several runs.\n *  - Remove the \"Professional Title\" placeholder paragraph.\n *  - Merge the two contact-info runs (phone/email + urls) into one run/line.\n *  - Merge each \"category heading\" + \"category detail\" paragraph pair in\n *    CORE COMPETENCIES into a single \"Label: detail\" paragraph.\n *  - Replace the single placeholder job with the full, real job history,\n *    using the Heading3 style for each job title line and plain \"\u2022\" bullet\n *    paragraphs (instead of the ListBullet style / glyphs) for the bullets.\n *  - Convert the KEY ACHIEVEMENTS heading run to use Heading3 and switch its\n *    checkmark bullets to plain \"\u2022\" paragraphs without the ListBullet style.\n *  - Widen the page margins.\n *\n * Implementation strategy: because so much of the body changes shape\n * (paragraphs merged, split, restyled, inserted, removed) the most robust\n * approach is to replace the whole body content in one shot with the target\n * WordprocessingML, then separately update the section page margins via the\n * dedicated pageSetup API (insertOoxml does not touch sectPr).\n */\n\nconst BODY_OOXML = \"<w:p><w:pPr><w:jc w:val=\\\"center\\\"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\\\"center\\\"/></w:pPr><w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading2\\\"/></w:pPr><w:r><w:t>PROFESSIONAL SUMMARY</w:t></w:r></w:p><w:p><w:r><w:t>Senior Software Engineer with 21 years building scalable geospatial data platforms, web applications, and distributed analytical systems. Expert in full-stack development with deep specialization in Apache Spark/Sedona for big data geospatial processing. Proven track record architecting multi-tenant SaaS platforms used by thousands of analysts, implementing ETL pipelines processing billions of geospatial records, and building production systems integrating ESRI, OSGeo, and SAFE FME technologies. Strong background in both enterprise consulting and startup environments, with experience leading engineering teams and delivering mission-critical geospatial applications.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading2\\\"/></w:pPr><w:r><w:t>CORE COMPETENCIES</w:t></w:r></w:p><w:p><w:r><w:t>Programming &amp; Development: Python: Django/GeoDjango, Flask, Pandas, PySpark, NumPy, SciKit-Learn \\u2022 JVM: Scala (Spark/Sedona), Java (GeoTools, enterprise applications), Groovy \\u2022 Web Technologies: JavaScript, React, d3.js, OpenLayers, jQuery, HTML/CSS \\u2022 Database Languages: SQL, T-SQL, PostgreSQL/PostGIS, Oracle, MySQL \\u2022 Statistical/Analysis: R, SPSS, NetLogo (agent-based modeling)</w:t></w:r></w:p><w:p><w:r><w:t>Big Data &amp; Geospatial Platforms: Apache Spark: PySpark, Spark SQL, Sedona (geospatial), distributed processing \\u2022 Geospatial Stack: PostGIS, ESRI ArcGIS, Quantum GIS, GRASS, OSGeo, SAFE FME \\u2022 Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Google Cloud, Microsoft Azure \\u2022 Data Engineering: ETL/ELT pipelines, dbt, Hadoop, Informatica, CDAP \\u2022 Databases: PostgreSQL/PostGIS, Oracle, MongoDB, Neo4j, MySQL</w:t></w:r></w:p><w:p><w:r><w:t>Software Architecture &amp; DevOps: Distributed Systems: Multi-tenant SaaS, microservices, API design, scalability \\u2022 Geospatial Applications: Spatial algorithms, boundary estimation, clustering analysis \\u2022 Web Applications: Full-stack development, RESTful APIs, real-time collaboration \\u2022 DevOps: Docker, Vagrant, CI/CD (GitLab, GitHub), Celery, Airflow, nginx \\u2022 Integration: Twilio API, WMS tile servers, CRM/DMP integration, OAuth</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading2\\\"/></w:pPr><w:r><w:t>PROFESSIONAL EXPERIENCE</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>PARTNER &amp; SENIOR SOFTWARE ENGINEER - Siege Analytics, Washington, DC | January 2014 \\u2013 Present</w:t></w:r></w:p><w:p><w:r><w:t>Geospatial Platform Architecture and Full-Stack Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Architected and engineered redistricting platform serving thousands of analysts with real-time collaborative editing, Census integration, and legal compliance analysis</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed boundary estimation microservice using incomplete data for boundary estimation without machine learning, processing geographies at national scale</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built scalable ETL pipelines using PySpark and Sedona processing billions of geospatial records with sub-hour latency requirements</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Implemented advanced spatial clustering algorithms achieving 88% improvement in analytical targeting efficacy for political applications</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Created fraud detection systems processing multi-terabyte campaign finance datasets with real-time alerting capabilities</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Led technical architecture decisions integrating ESRI, OSGeo, and SAFE FME technologies for Fortune 500 and political clients</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>PRINCIPAL SOFTWARE ENGINEER - Clarity and Rigour, Washington, DC | 2012 \\u2013 2014</w:t></w:r></w:p><w:p><w:r><w:t>Geospatial Solutions and Software Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed geospatial analysis frameworks and mapping applications for electoral research</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built custom visualization tools and interactive dashboards for client presentations</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Architected data processing pipelines for large-scale demographic and geographic datasets</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Created web-based mapping applications using JavaScript, OpenLayers, and PostGIS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>DIRECTOR OF DATA PRODUCTS - Helm, Washington, DC | 2010 \\u2013 2012</w:t></w:r></w:p><w:p><w:r><w:t>Data Product Development and Engineering Leadership</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Led development of data-driven solutions and platform architecture for political organizations</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Managed cross-functional engineering teams building campaign management and voter targeting systems</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Designed and implemented scalable data platforms using Python, Django, and PostgreSQL</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built RESTful APIs and microservices for campaign data integration</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>SENIOR SOFTWARE ENGINEER - GSD&amp;M, Austin, TX | 2008 \\u2013 2010</w:t></w:r></w:p><w:p><w:r><w:t>Campaign Technology and Analytics Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed voter targeting models and demographic analysis tools using Python and R</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built web applications for campaign data visualization and reporting</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Created data integration systems connecting multiple campaign data sources</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Implemented machine learning algorithms for voter behavior prediction</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 \\u2013 2006</w:t></w:r></w:p><w:p><w:r><w:t>Political Technology Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed software solutions for political campaigns and advocacy groups using PHP and JavaScript</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built web applications for voter engagement and campaign management</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Integrated third-party APIs and data sources for campaign tools</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Collaborated with political strategists to translate requirements into technical solutions</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \\u2013 2004</w:t></w:r></w:p><w:p><w:r><w:t>Nonprofit Technology Integration and Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed data management systems and web applications for social justice organizations</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built custom applications for community engagement and advocacy using PHP and MySQL</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Provided technical training and support to nonprofit staff</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Integrated technology solutions within organizational frameworks</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>PROGRAMMER - Lake Research Partners, Washington, DC | 2001 \\u2013 2002</w:t></w:r></w:p><w:p><w:r><w:t>Political Research and Data Analysis Tools</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed data analysis tools for political polling and research using Python and R</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built statistical models and data visualization tools for research presentations</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Created automated reporting systems for survey data analysis</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Supported senior researchers with technical analysis and data processing</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 \\u2013 2001</w:t></w:r></w:p><w:p><w:r><w:t>Political Field Operations and Data Management</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed data collection and management systems for field operations</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built databases and reporting tools for campaign field work</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Trained field staff on data collection protocols and quality control</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Analyzed field data to inform campaign strategy and research findings</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading2\\\"/></w:pPr><w:r><w:t>KEY ACHIEVEMENTS AND IMPACT</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Heading3\\\"/></w:pPr><w:r><w:t>Geospatial Platform Development</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Architected redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Built boundary estimation system achieving accurate geospatial results without machine learning using advanced PostGIS algorithms</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Developed econometric simulation platform with NetLogo multi-agent modeling and web interface</w:t></w:r></w:p><w:p><w:r><w:t>\\u2022 Created comprehensive survey platform managing complete research lifecycle with integrated geospatial market segmentation</w:t></w:r></w:p>\";\n\nconst pkg = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + BODY_OOXML + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nbody.insertOoxml(pkg, Word.InsertLocation.replace);\nawait context.sync();\n\n// Update page margins (twips -> points, 20 twips per point):\n//   top/bottom 864 -> 1440 (43.2pt -> 72pt), left/right 864 -> 1800 (43.2pt -> 90pt)\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst pageSetup = sections.items[0].pageSetup;\npageSetup.topMargin = 72;\npageSetup.bottomMargin = 72;\npageSetup.leftMargin = 90;\npageSetup.rightMargin = 90;\nawait context.sync();\n", "ps1": "# PowerShell / Word COM interop script implementing the resume rewrite:\n#  - Strip ad-hoc run colors/sizes from headings and several runs.\n#  - Remove the \"Professional Title\" placeholder paragraph.\n#  - Merge the two contact-info runs (phone/email + urls) into one run/line.\n#  - Merge each \"category heading\" + \"category detail\" paragraph pair in\n#    CORE COMPETENCIES into a single \"Label: detail\" paragraph.\n#  - Replace the single placeholder job with the full, real job history,\n#    using the Heading3 style for each job title line and plain \"bullet\n#    dot\" character paragraphs (instead of the ListBullet style / glyphs).\n#  - Convert the KEY ACHIEVEMENTS heading run to use Heading3 and switch its\n#    checkmark bullets to plain bullet paragraphs without the ListBullet style.\n#  - Widen the page margins.\n#\n# Implementation strategy: because so much of the body changes shape\n# (paragraphs merged, split, restyled, inserted, removed) the most robust\n# approach is to replace the whole body content in one shot via\n# Range.InsertXML (WordprocessingML wrapped in the xmlPackage part format),\n# then separately update the section page margins via PageSetup (InsertXML\n# only touches the body content, not sectPr).\n\n$d = $word.ActiveDocument\n\n$bodyOoxml = @'\n<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>PROFESSIONAL SUMMARY</w:t></w:r></w:p><w:p><w:r><w:t>Senior Software Engineer with 21 years building scalable geospatial data platforms, web applications, and distributed analytical systems. Expert in full-stack development with deep specialization in Apache Spark/Sedona for big data geospatial processing. Proven track record architecting multi-tenant SaaS platforms used by thousands of analysts, implementing ETL pipelines processing billions of geospatial records, and building production systems integrating ESRI, OSGeo, and SAFE FME technologies. Strong background in both enterprise consulting and startup environments, with experience leading engineering teams and delivering mission-critical geospatial applications.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>CORE COMPETENCIES</w:t></w:r></w:p><w:p><w:r><w:t>Programming &amp; Development: Python: Django/GeoDjango, Flask, Pandas, PySpark, NumPy, SciKit-Learn \u2022 JVM: Scala (Spark/Sedona), Java (GeoTools, enterprise applications), Groovy \u2022 Web Technologies: JavaScript, React, d3.js, OpenLayers, jQuery, HTML/CSS \u2022 Database Languages: SQL, T-SQL, PostgreSQL/PostGIS, Oracle, MySQL \u2022 Statistical/Analysis: R, SPSS, NetLogo (agent-based modeling)</w:t></w:r></w:p><w:p><w:r><w:t>Big Data &amp; Geospatial Platforms: Apache Spark: PySpark, Spark SQL, Sedona (geospatial), distributed processing \u2022 Geospatial Stack: PostGIS, ESRI ArcGIS, Quantum GIS, GRASS, OSGeo, SAFE FME \u2022 Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Google Cloud, Microsoft Azure \u2022 Data Engineering: ETL/ELT pipelines, dbt, Hadoop, Informatica, CDAP \u2022 Databases: PostgreSQL/PostGIS, Oracle, MongoDB, Neo4j, MySQL</w:t></w:r></w:p><w:p><w:r><w:t>Software Architecture &amp; DevOps: Distributed Systems: Multi-tenant SaaS, microservices, API design, scalability \u2022 Geospatial Applications: Spatial algorithms, boundary estimation, clustering analysis \u2022 Web Applications: Full-stack development, RESTful APIs, real-time collaboration \u2022 DevOps: Docker, Vagrant, CI/CD (GitLab, GitHub), Celery, Airflow, nginx \u2022 Integration: Twilio API, WMS tile servers, CRM/DMP integration, OAuth</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>PROFESSIONAL EXPERIENCE</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>PARTNER &amp; SENIOR SOFTWARE ENGINEER - Siege Analytics, Washington, DC | January 2014 \u2013 Present</w:t></w:r></w:p><w:p><w:r><w:t>Geospatial Platform Architecture and Full-Stack Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Architected and engineered redistricting platform serving thousands of analysts with real-time collaborative editing, Census integration, and legal compliance analysis</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed boundary estimation microservice using incomplete data for boundary estimation without machine learning, processing geographies at national scale</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built scalable ETL pipelines using PySpark and Sedona processing billions of geospatial records with sub-hour latency requirements</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Implemented advanced spatial clustering algorithms achieving 88% improvement in analytical targeting efficacy for political applications</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Created fraud detection systems processing multi-terabyte campaign finance datasets with real-time alerting capabilities</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Led technical architecture decisions integrating ESRI, OSGeo, and SAFE FME technologies for Fortune 500 and political clients</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>PRINCIPAL SOFTWARE ENGINEER - Clarity and Rigour, Washington, DC | 2012 \u2013 2014</w:t></w:r></w:p><w:p><w:r><w:t>Geospatial Solutions and Software Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed geospatial analysis frameworks and mapping applications for electoral research</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built custom visualization tools and interactive dashboards for client presentations</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Architected data processing pipelines for large-scale demographic and geographic datasets</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Created web-based mapping applications using JavaScript, OpenLayers, and PostGIS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>DIRECTOR OF DATA PRODUCTS - Helm, Washington, DC | 2010 \u2013 2012</w:t></w:r></w:p><w:p><w:r><w:t>Data Product Development and Engineering Leadership</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Led development of data-driven solutions and platform architecture for political organizations</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Managed cross-functional engineering teams building campaign management and voter targeting systems</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Designed and implemented scalable data platforms using Python, Django, and PostgreSQL</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built RESTful APIs and microservices for campaign data integration</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>SENIOR SOFTWARE ENGINEER - GSD&amp;M, Austin, TX | 2008 \u2013 2010</w:t></w:r></w:p><w:p><w:r><w:t>Campaign Technology and Analytics Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed voter targeting models and demographic analysis tools using Python and R</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built web applications for campaign data visualization and reporting</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Created data integration systems connecting multiple campaign data sources</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Implemented machine learning algorithms for voter behavior prediction</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 \u2013 2006</w:t></w:r></w:p><w:p><w:r><w:t>Political Technology Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed software solutions for political campaigns and advocacy groups using PHP and JavaScript</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built web applications for voter engagement and campaign management</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Integrated third-party APIs and data sources for campaign tools</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Collaborated with political strategists to translate requirements into technical solutions</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \u2013 2004</w:t></w:r></w:p><w:p><w:r><w:t>Nonprofit Technology Integration and Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed data management systems and web applications for social justice organizations</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built custom applications for community engagement and advocacy using PHP and MySQL</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Provided technical training and support to nonprofit staff</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Integrated technology solutions within organizational frameworks</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>PROGRAMMER - Lake Research Partners, Washington, DC | 2001 \u2013 2002</w:t></w:r></w:p><w:p><w:r><w:t>Political Research and Data Analysis Tools</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed data analysis tools for political polling and research using Python and R</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built statistical models and data visualization tools for research presentations</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Created automated reporting systems for survey data analysis</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Supported senior researchers with technical analysis and data processing</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 \u2013 2001</w:t></w:r></w:p><w:p><w:r><w:t>Political Field Operations and Data Management</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed data collection and management systems for field operations</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built databases and reporting tools for campaign field work</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Trained field staff on data collection protocols and quality control</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Analyzed field data to inform campaign strategy and research findings</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>KEY ACHIEVEMENTS AND IMPACT</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:r><w:t>Geospatial Platform Development</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Architected redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Built boundary estimation system achieving accurate geospatial results without machine learning using advanced PostGIS algorithms</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Developed econometric simulation platform with NetLogo multi-agent modeling and web interface</w:t></w:r></w:p><w:p><w:r><w:t>\u2022 Created comprehensive survey platform managing complete research lifecycle with integrated geospatial market segmentation</w:t></w:r></w:p>\n'@\n\n$pkg = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>' + $bodyOoxml + '</w:body>' + `\n  '</w:document>' + `\n  '</pkg:xmlData>' + `\n  '</pkg:part>' + `\n  '</pkg:package>'\n\n$d.Content.InsertXML($pkg)\n\n# Update page margins (twips -> points, 20 twips per point):\n#   top/bottom 864 -> 1440 (43.2pt -> 72pt), left/right 864 -> 1800 (43.2pt -> 90pt)\n$d.PageSetup.TopMargin = 72\n$d.PageSetup.BottomMargin = 72\n$d.PageSetup.LeftMargin = 90\n$d.PageSetup.RightMargin = 90\n"}
